# Applies the "Updated symbol list" GitHub Actions commit to the crypto price sheet.
# Each cell holds an inline string (even the numeric-looking price column), so every
# write is forced to Text via a leading apostrophe and the cell style is reset to
# "Normal" afterwards so Excel does not silently coerce the value to a Number.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $text) {
    $range = $ws.Range($cellRef)
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

Set-TextValue "D2" '241.77'
Set-TextValue "D3" '21.80'
Set-TextValue "D4" '5.391'
Set-TextValue "D5" '0.05678'
Set-TextValue "D7" '6.279'
Set-TextValue "D8" '0.8070'
Set-TextValue "D9" '0.9859'
Set-TextValue "B10" 'One'
Set-TextValue "C10" 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
Set-TextValue "D10" '0.01067'
Set-TextValue "E10" '9OneONEBestin24h'
Set-TextValue "B11" 'WazirX'
Set-TextValue "C11" 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
Set-TextValue "D11" '0.1419'
Set-TextValue "E11" '10WazirXWRX'
Set-TextValue "B12" 'MandalaExchangeToken'
Set-TextValue "C12" 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
Set-TextValue "D12" '0.07280'
Set-TextValue "E12" '11MandalaExchangeTokenMDX'
Set-TextValue "B13" 'LiechtensteinCryptoassetsExchange'
Set-TextValue "C13" 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
Set-TextValue "D13" '0.03039'
Set-TextValue "E13" '12LiechtensteinCryptoassetsExchangeLCX'
Set-TextValue "B14" 'BitrueCoin'
Set-TextValue "C14" 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
Set-TextValue "D14" '0.03122'
Set-TextValue "E14" '13BitrueCoinBTR'
Set-TextValue "B15" 'ProBitToken'
Set-TextValue "C15" 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'
Set-TextValue "D15" '0.1310'
Set-TextValue "E15" '14ProBitTokenPROB'
Set-TextValue "B16" 'BitMartToken'
Set-TextValue "C16" 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
Set-TextValue "D16" '0.09352'
Set-TextValue "E16" '15BitMartTokenBMX'
Set-TextValue "B17" 'MCDex'
Set-TextValue "C17" 'https://coinranking.com/coin/3nMM61qeg+mcdex-mcb'
Set-TextValue "D17" '3.916'
Set-TextValue "E17" '16MCDexMCB'
Set-TextValue "B18" 'BitForexToken'
Set-TextValue "C18" 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
Set-TextValue "D18" '0.001597'
Set-TextValue "E18" '17BitForexTokenBF'
Set-TextValue "B19" 'CoinExToken'
Set-TextValue "C19" 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'
Set-TextValue "D19" '0.04789'
Set-TextValue "E19" '18CoinExTokenCET'
Set-TextValue "B20" 'TigerCash'
Set-TextValue "C20" 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
Set-TextValue "D20" '0.006293'
Set-TextValue "E20" '19TigerCashTCH'
Set-TextValue "B21" 'HotbitToken'
Set-TextValue "C21" 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'
Set-TextValue "D21" '0.004056'
Set-TextValue "E21" '20HotbitTokenHTB'
Set-TextValue "B22" 'BitKan'
Set-TextValue "C22" 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'
Set-TextValue "D22" '0.0009919'
Set-TextValue "E22" '21BitKanKAN'
Set-TextValue "B23" 'NitroEx'
Set-TextValue "C23" 'https://coinranking.com/coin/8oiZw6gwYhC+nitroex-ntx'
Set-TextValue "D23" '0.0001500'
Set-TextValue "E23" '22NitroExNTX'
Set-TextValue "B24" 'LEO'
Set-TextValue "C24" 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
Set-TextValue "D24" '3.739'
Set-TextValue "E24" '23LEOLEO'
Set-TextValue "B25" 'BTSEToken'
Set-TextValue "C25" 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
Set-TextValue "D25" '2.154'
Set-TextValue "E25" '24BTSETokenBTSE'
Set-TextValue "B26" 'BitpandaEcosystemToken'
Set-TextValue "C26" 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
Set-TextValue "D26" '0.3259'
Set-TextValue "E26" '25BitpandaEcosystemTokenBEST'
Set-TextValue "D27" '0.0004000'
Set-TextValue "D40" '0.03804'
Set-TextValue "D41" '0.006669'
Set-TextValue "D42" '0.1045'
Set-TextValue "D44" '0.006813'
Set-TextValue "D45" '0.00005613'
Set-TextValue "D47" '0.3901'
Set-TextValue "E47" '46CoinbaseStockTokenCOIN'
Set-TextValue "D49" '0.00002100'
